$wb = $excel.ActiveWorkbook

# --- 1. Text change: "Ready for handoff" -> "In Translation" ---
# This shared string is used by the per-language status cells on all
# three sheets (Overview!E2:F3 and the Status column C2:C3 on the
# zh-cn / de-de detail sheets).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- 2. Column width changes ---
# Overview columns E and F (zh-cn / de-de status columns) get narrower,
# and the Status column (C) on the zh-cn / de-de detail sheets matches.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
